# config thông tin trong các sheet
# Updates the "last_edited_time" (column D) values for the Notion rows to
# reflect the new edit timestamps recorded by Notion on 2024-08-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value        = "2024-08-31T05:43:00.000Z"
$ws.Range("D3:D6").Value     = "2024-08-31T05:39:00.000Z"
$ws.Range("D7:D61").Value    = "2024-08-31T05:40:00.000Z"
$ws.Range("D62:D116").Value  = "2024-08-31T05:41:00.000Z"
$ws.Range("D117:D157").Value = "2024-08-31T05:42:00.000Z"
